# Weekly data update (Fruta / hortaliza, semanal):
# A new price-report row for "Pina" (Segunda quality, $/caja 14 unidades) is
# inserted at row 385, shifting every subsequent record down by one row
# through row 472. Rewritten here as a full row-by-row value rewrite of
# A385:T472 (equivalent end state to inserting a row and filling it in).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45204, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 25000, 25000, 25000, '$/caja 14 unidades', 'Ecuador', 1786, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(385, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44414, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 160, 22000, 22000, 22000, '$/caja 14 unidades', 'Ecuador', 1571, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(386, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44571, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 50, 19000, 20000, 19500, '$/caja 16 unidades', 'Ecuador', 1219, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(387, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44258, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 20, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(388, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44238, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 120, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(389, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44925, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 300, 19000, 20000, 19500, '$/caja 14 unidades', 'Ecuador', 1393, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(390, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44446, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 160, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(391, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44340, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 120, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(392, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44246, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 160, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(393, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44494, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 180, 25000, 25500, 25250, '$/caja 14 unidades', 'Ecuador', 1804, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(394, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44817, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(395, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44323, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 200, 18000, 19000, 18500, '$/caja 14 unidades', 'Ecuador', 1321, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(396, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44838, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 23500, 24000, 23750, '$/caja 12 unidades', 'Ecuador', 1979, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(397, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44551, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 300, 20000, 21000, 20500, '$/caja 16 unidades', 'Ecuador', 1281, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(398, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44830, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 60, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(399, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44678, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 30, 17000, 18000, 17500, '$/caja 16 unidades', 'Ecuador', 1094, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(400, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45125, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 26000, 26000, 26000, '$/caja 14 unidades', 'Ecuador', 1857, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(401, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45124, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 30, 27000, 27000, 27000, '$/caja 12 unidades', 'Ecuador', 2250, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(402, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45135, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 150, 22000, 22000, 22000, '$/caja 14 unidades', 'Ecuador', 1571, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(403, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44578, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 120, 19000, 20000, 19500, '$/caja 16 unidades', 'Ecuador', 1219, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(404, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44799, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 100, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(405, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44704, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 80, 22000, 22000, 22000, '$/caja 16 unidades', 'Ecuador', 1375, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(406, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44694, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 400, 22000, 23000, 22500, '$/caja 16 unidades', 'Ecuador', 1406, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(407, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44498, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 360, 20000, 21000, 20500, '$/caja 14 unidades', 'Ecuador', 1464, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(408, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44859, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 160, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(409, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44845, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(410, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44363, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 80, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(411, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44406, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 80, 24000, 24000, 24000, '$/caja 14 unidades', 'Ecuador', 1714, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(412, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44747, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 160, 22000, 22000, 22000, '$/caja 16 unidades', 'Ecuador', 1375, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(413, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44924, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 120, 19000, 20000, 19500, '$/caja 14 unidades', 'Ecuador', 1393, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(414, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45055, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 20000, 21000, 20500, '$/caja 12 unidades', 'Ecuador', 1708, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(415, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44315, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 120, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(416, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45044, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 20000, 21000, 20500, '$/caja 12 unidades', 'Ecuador', 1708, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(417, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45044, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 200, 18000, 19000, 18500, '$/caja 14 unidades', 'Ecuador', 1321, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(418, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44518, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 120, 21000, 22000, 21500, '$/caja 16 unidades', 'Ecuador', 1344, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(419, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44701, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 300, 21000, 22000, 21500, '$/caja 16 unidades', 'Ecuador', 1344, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(420, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44824, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 20, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(421, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44981, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 26000, 27000, 26500, '$/caja 12 unidades', 'Ecuador', 2208, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(422, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45190, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 25000, 25000, 25000, '$/caja 14 unidades', 'Ecuador', 1786, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(423, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44264, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 200, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(424, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45043, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 120, 20000, 21000, 20500, '$/caja 12 unidades', 'Ecuador', 1708, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(425, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44396, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 20000, 20000, 20000, '$/caja 14 unidades', 'Ecuador', 1429, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(426, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44449, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 160, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(427, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44249, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 17500, 18000, 17750, '$/caja 14 unidades', 'Ecuador', 1268, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(428, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45071, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 80, 19000, 20000, 19500, '$/caja 12 unidades', 'Ecuador', 1625, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(429, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44377, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 80, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(430, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44438, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 80, 21000, 21000, 21000, '$/caja 12 unidades', 'Ecuador', 1750, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(431, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44536, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 120, 20000, 21000, 20500, '$/caja 16 unidades', 'Ecuador', 1281, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(432, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45033, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 2000, 21000, 11500, '$/caja 12 unidades', 'Ecuador', 958, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(433, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44558, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 200, 19000, 20000, 19500, '$/caja 16 unidades', 'Ecuador', 1219, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(434, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44364, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 120, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(435, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44624, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 300, 19000, 20000, 19500, '$/caja 16 unidades', 'Ecuador', 1219, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(436, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44358, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 160, 17000, 18000, 17500, '$/caja 14 unidades', 'Colombia', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(437, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45134, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 80, 22000, 22000, 22000, '$/caja 14 unidades', 'Ecuador', 1571, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(438, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44792, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 120, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(439, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44757, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 140, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(440, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44813, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(441, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44419, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 40, 21000, 21000, 21000, '$/caja 14 unidades', 'Ecuador', 1500, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(442, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44918, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 300, 22000, 23000, 22500, '$/caja 14 unidades', 'Ecuador', 1607, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(443, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45077, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 40, 18000, 19000, 18500, '$/caja 12 unidades', 'Ecuador', 1542, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(444, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44335, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 60, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(445, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44434, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 80, 22000, 22000, 22000, '$/caja 12 unidades', 'Ecuador', 1833, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(446, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44798, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 100, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(447, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45005, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 21000, 22000, 21500, '$/caja 12 unidades', 'Ecuador', 1792, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(448, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45175, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 40, 25000, 25000, 25000, '$/caja 14 unidades', 'Ecuador', 1786, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(449, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45194, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 25000, 25000, 25000, '$/caja 14 unidades', 'Ecuador', 1786, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(450, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44511, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 90, 22000, 23000, 22500, '$/caja 14 unidades', 'Ecuador', 1607, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(451, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44579, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 200, 19000, 20000, 19500, '$/caja 16 unidades', 'Ecuador', 1219, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(452, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44767, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 20, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(453, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44973, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 70, 27000, 28000, 27500, '$/caja 12 unidades', 'Ecuador', 2292, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(454, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44665, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 240, 15000, 16000, 15500, '$/caja 16 unidades', 'Ecuador', 969, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(455, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44544, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 200, 20000, 21000, 20500, '$/caja 16 unidades', 'Ecuador', 1281, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(456, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44825, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 20, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(457, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44432, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 140, 22000, 22000, 22000, '$/caja 12 unidades', 'Ecuador', 1833, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(458, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44939, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 300, 21000, 22000, 21500, '$/caja 14 unidades', 'Ecuador', 1536, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(459, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44803, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 100, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(460, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45176, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 100, 25000, 25000, 25000, '$/caja 14 unidades', 'Ecuador', 1786, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(461, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44960, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 200, 23000, 24000, 23500, '$/caja 12 unidades', 'Ecuador', 1958, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(462, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44832, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 20, 23000, 23500, 23250, '$/caja 12 unidades', 'Ecuador', 1938, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(463, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 45007, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 80, 20000, 21000, 20500, '$/caja 12 unidades', 'Ecuador', 1708, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(464, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44594, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 80, 17000, 18000, 17500, '$/caja 16 unidades', 'Ecuador', 1094, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(465, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44586, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 200, 19000, 19000, 19000, '$/caja 16 unidades', 'Ecuador', 1188, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(466, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44567, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 120, 19000, 20000, 19500, '$/caja 16 unidades', 'Ecuador', 1219, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(467, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44270, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 40, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(468, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44243, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 160, 17000, 18000, 17500, '$/caja 14 unidades', 'Ecuador', 1250, 14)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(469, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44809, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 100, 23000, 23000, 23000, '$/caja 12 unidades', 'Ecuador', 1917, 12)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(470, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44698, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 200, 22000, 23000, 22500, '$/caja 16 unidades', 'Ecuador', 1406, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(471, $col).Value = $rowData[$col - 1] }
$rowData = @(4, 'Feria Lagunitas de Puerto Montt', 'Los Lagos', 44552, 10, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 60, 20000, 21000, 20500, '$/caja 16 unidades', 'Ecuador', 1281, 16)
for ($col = 1; $col -le 20; $col++) { $ws.Cells.Item(472, $col).Value = $rowData[$col - 1] }

# Preserve the date number-format on column D (newly created row 472 otherwise defaults to General).
$ws.Range('D385:D472').NumberFormat = 'YYYY-MM-DD HH:MM:SS'

$ws.Range('A1').Select()
